# Auto-generated Excel COM-interop script to apply recomputed pipeline values
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Step1_Data ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 4).Value2 = 0.1331171958662623
$ws.Cells.Item(2, 5).Value2 = 0.005823602207493439
$ws.Cells.Item(2, 6).Value2 = 0.146611740865868
$ws.Cells.Item(2, 12).Value2 = 0.03594461592664443
$ws.Cells.Item(2, 13).Value2 = 0.0184067265993822
$ws.Cells.Item(2, 14).Value2 = 0.1421412447461528
$ws.Cells.Item(2, 17).Value2 = 0.01349153145915159
$ws.Cells.Item(2, 18).Value2 = 0.04773208950185645
$ws.Cells.Item(2, 19).Value2 = 0.1401666395175499
$ws.Cells.Item(2, 20).Value2 = 0.01226389577420914
$ws.Cells.Item(2, 21).Value2 = 0.03823009365066198
$ws.Cells.Item(2, 22).Value2 = 0.06638765525501385
$ws.Cells.Item(2, 23).Value2 = 0.02137546217872988
$ws.Cells.Item(2, 24).Value2 = 0.05402735103370124
$ws.Cells.Item(2, 26).Value2 = 0.1242801554173228
$ws.Cells.Item(3, 4).Value2 = 0.1322594850399036
$ws.Cells.Item(3, 6).Value2 = 0.1406990439050197
$ws.Cells.Item(3, 7).Value2 = 0.01019179423398317
$ws.Cells.Item(3, 8).Value2 = 0.002963224719556238
$ws.Cells.Item(3, 12).Value2 = 0.03602714371673011
$ws.Cells.Item(3, 14).Value2 = 0.1602398646078278
$ws.Cells.Item(3, 15).Value2 = 0.00587677650074706
$ws.Cells.Item(3, 17).Value2 = 0.00781303421535179
$ws.Cells.Item(3, 19).Value2 = 0.1908142222560892
$ws.Cells.Item(3, 21).Value2 = 0.07888585070747095
$ws.Cells.Item(3, 22).Value2 = 0.004690765684915261
$ws.Cells.Item(3, 23).Value2 = 0.09644656129453208
$ws.Cells.Item(3, 24).Value2 = 0.01762989566923217
$ws.Cells.Item(3, 25).Value2 = 0.004309507020818002
$ws.Cells.Item(3, 26).Value2 = 0.1046637898627124
$ws.Cells.Item(3, 27).Value2 = 0.006489040565110397
$ws.Cells.Item(4, 4).Value2 = 0.06663486524591113
$ws.Cells.Item(4, 6).Value2 = 0.1727373926859226
$ws.Cells.Item(4, 7).Value2 = 0.02557076284951977
$ws.Cells.Item(4, 11).Value2 = 0.002449898481586343
$ws.Cells.Item(4, 12).Value2 = 0.03385597157904392
$ws.Cells.Item(4, 14).Value2 = 0.1918617609668229
$ws.Cells.Item(4, 15).Value2 = 0.0322600932851816
$ws.Cells.Item(4, 17).Value2 = 0.0201666074228898
$ws.Cells.Item(4, 19).Value2 = 0.1551396586254308
$ws.Cells.Item(4, 21).Value2 = 0.0949170166081893
$ws.Cells.Item(4, 23).Value2 = 0.07447265087109206
$ws.Cells.Item(4, 24).Value2 = 0.02090400901499286
$ws.Cells.Item(4, 25).Value2 = 0.01274682718887407
$ws.Cells.Item(4, 26).Value2 = 0.08578713501260242
$ws.Cells.Item(4, 27).Value2 = 0.01012534439722046
$ws.Cells.Item(4, 28).Value2 = 0.0003700057647198986
$ws.Cells.Item(5, 4).Value2 = 0.02973020990429905
$ws.Cells.Item(5, 5).Value2 = 0.02537520941719335
$ws.Cells.Item(5, 6).Value2 = 0.1431374812087675
$ws.Cells.Item(5, 7).Value2 = 0.03483576522052386
$ws.Cells.Item(5, 11).Value2 = 0.02147117435759112
$ws.Cells.Item(5, 13).Value2 = 0.03323073099509084
$ws.Cells.Item(5, 14).Value2 = 0.1367198572579337
$ws.Cells.Item(5, 15).Value2 = 0.08794296524575976
$ws.Cells.Item(5, 19).Value2 = 0.1376866995076497
$ws.Cells.Item(5, 20).Value2 = 0.02930544316053927
$ws.Cells.Item(5, 21).Value2 = 0.07749576014494995
$ws.Cells.Item(5, 23).Value2 = 0.09873724827680609
$ws.Cells.Item(5, 24).Value2 = 0.001058114437597048
$ws.Cells.Item(5, 25).Value2 = 0.02616611974039606
$ws.Cells.Item(5, 26).Value2 = 0.07762778566850162
$ws.Cells.Item(5, 27).Value2 = 0.03902588647276124
$ws.Cells.Item(5, 28).Value2 = 0.0004535489836398096
$ws.Cells.Item(6, 4).Value2 = 0.09315088171296432
$ws.Cells.Item(6, 6).Value2 = 0.1548668062577585
$ws.Cells.Item(6, 7).Value2 = 0.02824288946482174
$ws.Cells.Item(6, 12).Value2 = 0.02860327889003974
$ws.Cells.Item(6, 14).Value2 = 0.1295439262458098
$ws.Cells.Item(6, 15).Value2 = 0.04095276189164863
$ws.Cells.Item(6, 19).Value2 = 0.1631654012846804
$ws.Cells.Item(6, 21).Value2 = 0.1118037266411649
$ws.Cells.Item(6, 23).Value2 = 0.1166085519814386
$ws.Cells.Item(6, 24).Value2 = 0.005560601891444221
$ws.Cells.Item(6, 25).Value2 = 0.03360123739185195
$ws.Cells.Item(6, 26).Value2 = 0.06844750439368316
$ws.Cells.Item(6, 27).Value2 = 0.02545243195269409

# --- Sheet 2: Step2_Sj ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 4).Value2 = 0.1331171958662623
$ws.Cells.Item(2, 5).Value2 = 0.1389407980737557
$ws.Cells.Item(2, 6).Value2 = 0.2855525389396237
$ws.Cells.Item(2, 7).Value2 = 0.2855525389396237
$ws.Cells.Item(2, 8).Value2 = 0.2855525389396237
$ws.Cells.Item(2, 9).Value2 = 0.2855525389396237
$ws.Cells.Item(2, 10).Value2 = 0.2855525389396237
$ws.Cells.Item(2, 11).Value2 = 0.2855525389396237
$ws.Cells.Item(2, 12).Value2 = 0.3214971548662682
$ws.Cells.Item(2, 13).Value2 = 0.3399038814656504
$ws.Cells.Item(2, 14).Value2 = 0.4820451262118032
$ws.Cells.Item(2, 15).Value2 = 0.4820451262118032
$ws.Cells.Item(2, 16).Value2 = 0.4820451262118032
$ws.Cells.Item(2, 17).Value2 = 0.4955366576709548
$ws.Cells.Item(2, 18).Value2 = 0.5432687471728113
$ws.Cells.Item(2, 19).Value2 = 0.6834353866903612
$ws.Cells.Item(2, 20).Value2 = 0.6956992824645704
$ws.Cells.Item(2, 21).Value2 = 0.7339293761152323
$ws.Cells.Item(2, 22).Value2 = 0.8003170313702462
$ws.Cells.Item(2, 23).Value2 = 0.8216924935489761
$ws.Cells.Item(2, 24).Value2 = 0.8757198445826774
$ws.Cells.Item(2, 25).Value2 = 0.8757198445826774
$ws.Cells.Item(3, 4).Value2 = 0.1322594850399036
$ws.Cells.Item(3, 5).Value2 = 0.1322594850399036
$ws.Cells.Item(3, 6).Value2 = 0.2729585289449233
$ws.Cells.Item(3, 7).Value2 = 0.2831503231789065
$ws.Cells.Item(3, 8).Value2 = 0.2861135478984627
$ws.Cells.Item(3, 9).Value2 = 0.2861135478984627
$ws.Cells.Item(3, 10).Value2 = 0.2861135478984627
$ws.Cells.Item(3, 11).Value2 = 0.2861135478984627
$ws.Cells.Item(3, 12).Value2 = 0.3221406916151928
$ws.Cells.Item(3, 13).Value2 = 0.3221406916151928
$ws.Cells.Item(3, 14).Value2 = 0.4823805562230206
$ws.Cells.Item(3, 15).Value2 = 0.4882573327237677
$ws.Cells.Item(3, 16).Value2 = 0.4882573327237677
$ws.Cells.Item(3, 17).Value2 = 0.4960703669391195
$ws.Cells.Item(3, 18).Value2 = 0.4960703669391195
$ws.Cells.Item(3, 19).Value2 = 0.6868845891952087
$ws.Cells.Item(3, 20).Value2 = 0.6868845891952087
$ws.Cells.Item(3, 21).Value2 = 0.7657704399026797
$ws.Cells.Item(3, 22).Value2 = 0.7704612055875949
$ws.Cells.Item(3, 23).Value2 = 0.866907766882127
$ws.Cells.Item(3, 24).Value2 = 0.8845376625513591
$ws.Cells.Item(3, 25).Value2 = 0.8888471695721771
$ws.Cells.Item(3, 26).Value2 = 0.9935109594348895
$ws.Cells.Item(3, 27).Value2 = 0.9999999999999999
$ws.Cells.Item(3, 28).Value2 = 0.9999999999999999
$ws.Cells.Item(3, 29).Value2 = 0.9999999999999999
$ws.Cells.Item(3, 30).Value2 = 0.9999999999999999
$ws.Cells.Item(3, 31).Value2 = 0.9999999999999999
$ws.Cells.Item(3, 32).Value2 = 0.9999999999999999
$ws.Cells.Item(3, 33).Value2 = 0.9999999999999999
$ws.Cells.Item(3, 34).Value2 = 0.9999999999999999
$ws.Cells.Item(3, 35).Value2 = 0.9999999999999999
$ws.Cells.Item(4, 4).Value2 = 0.06663486524591113
$ws.Cells.Item(4, 5).Value2 = 0.06663486524591113
$ws.Cells.Item(4, 6).Value2 = 0.2393722579318338
$ws.Cells.Item(4, 7).Value2 = 0.2649430207813535
$ws.Cells.Item(4, 8).Value2 = 0.2649430207813535
$ws.Cells.Item(4, 9).Value2 = 0.2649430207813535
$ws.Cells.Item(4, 10).Value2 = 0.2649430207813535
$ws.Cells.Item(4, 11).Value2 = 0.2673929192629398
$ws.Cells.Item(4, 12).Value2 = 0.3012488908419838
$ws.Cells.Item(4, 13).Value2 = 0.3012488908419838
$ws.Cells.Item(4, 14).Value2 = 0.4931106518088066
$ws.Cells.Item(4, 15).Value2 = 0.5253707450939883
$ws.Cells.Item(4, 16).Value2 = 0.5253707450939883
$ws.Cells.Item(4, 17).Value2 = 0.5455373525168781
$ws.Cells.Item(4, 18).Value2 = 0.5455373525168781
$ws.Cells.Item(4, 19).Value2 = 0.7006770111423088
$ws.Cells.Item(4, 20).Value2 = 0.7006770111423088
$ws.Cells.Item(4, 21).Value2 = 0.7955940277504981
$ws.Cells.Item(4, 22).Value2 = 0.7955940277504981
$ws.Cells.Item(4, 23).Value2 = 0.8700666786215903
$ws.Cells.Item(4, 24).Value2 = 0.8909706876365832
$ws.Cells.Item(4, 25).Value2 = 0.9037175148254573
$ws.Cells.Item(4, 26).Value2 = 0.9895046498380596
$ws.Cells.Item(4, 27).Value2 = 0.9996299942352801
$ws.Cells.Item(5, 4).Value2 = 0.02973020990429905
$ws.Cells.Item(5, 5).Value2 = 0.0551054193214924
$ws.Cells.Item(5, 6).Value2 = 0.1982429005302599
$ws.Cells.Item(5, 7).Value2 = 0.2330786657507838
$ws.Cells.Item(5, 8).Value2 = 0.2330786657507838
$ws.Cells.Item(5, 9).Value2 = 0.2330786657507838
$ws.Cells.Item(5, 10).Value2 = 0.2330786657507838
$ws.Cells.Item(5, 11).Value2 = 0.2545498401083749
$ws.Cells.Item(5, 12).Value2 = 0.2545498401083749
$ws.Cells.Item(5, 13).Value2 = 0.2877805711034658
$ws.Cells.Item(5, 14).Value2 = 0.4245004283613995
$ws.Cells.Item(5, 15).Value2 = 0.5124433936071593
$ws.Cells.Item(5, 16).Value2 = 0.5124433936071593
$ws.Cells.Item(5, 17).Value2 = 0.5124433936071593
$ws.Cells.Item(5, 18).Value2 = 0.5124433936071593
$ws.Cells.Item(5, 19).Value2 = 0.650130093114809
$ws.Cells.Item(5, 20).Value2 = 0.6794355362753483
$ws.Cells.Item(5, 21).Value2 = 0.7569312964202982
$ws.Cells.Item(5, 22).Value2 = 0.7569312964202982
$ws.Cells.Item(5, 23).Value2 = 0.8556685446971043
$ws.Cells.Item(5, 24).Value2 = 0.8567266591347014
$ws.Cells.Item(5, 25).Value2 = 0.8828927788750974
$ws.Cells.Item(5, 26).Value2 = 0.960520564543599
$ws.Cells.Item(5, 27).Value2 = 0.9995464510163603
$ws.Cells.Item(6, 4).Value2 = 0.09315088171296432
$ws.Cells.Item(6, 5).Value2 = 0.09315088171296432
$ws.Cells.Item(6, 6).Value2 = 0.2480176879707228
$ws.Cells.Item(6, 7).Value2 = 0.2762605774355446
$ws.Cells.Item(6, 8).Value2 = 0.2762605774355446
$ws.Cells.Item(6, 9).Value2 = 0.2762605774355446
$ws.Cells.Item(6, 10).Value2 = 0.2762605774355446
$ws.Cells.Item(6, 11).Value2 = 0.2762605774355446
$ws.Cells.Item(6, 12).Value2 = 0.3048638563255843
$ws.Cells.Item(6, 13).Value2 = 0.3048638563255843
$ws.Cells.Item(6, 14).Value2 = 0.4344077825713941
$ws.Cells.Item(6, 15).Value2 = 0.4753605444630428
$ws.Cells.Item(6, 16).Value2 = 0.4753605444630428
$ws.Cells.Item(6, 17).Value2 = 0.4753605444630428
$ws.Cells.Item(6, 18).Value2 = 0.4753605444630428
$ws.Cells.Item(6, 19).Value2 = 0.6385259457477233
$ws.Cells.Item(6, 20).Value2 = 0.6385259457477233
$ws.Cells.Item(6, 21).Value2 = 0.7503296723888881
$ws.Cells.Item(6, 22).Value2 = 0.7503296723888881
$ws.Cells.Item(6, 23).Value2 = 0.8669382243703267
$ws.Cells.Item(6, 24).Value2 = 0.872498826261771
$ws.Cells.Item(6, 25).Value2 = 0.9061000636536229
$ws.Cells.Item(6, 26).Value2 = 0.9745475680473061
$ws.Cells.Item(6, 27).Value2 = 1
$ws.Cells.Item(6, 28).Value2 = 1
$ws.Cells.Item(6, 29).Value2 = 1
$ws.Cells.Item(6, 30).Value2 = 1
$ws.Cells.Item(6, 31).Value2 = 1
$ws.Cells.Item(6, 32).Value2 = 1
$ws.Cells.Item(6, 33).Value2 = 1
$ws.Cells.Item(6, 34).Value2 = 1
$ws.Cells.Item(6, 35).Value2 = 1

# --- Sheet 3: Step3_DataPts_0.5 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 4).Value2 = 17
$ws.Cells.Item(2, 6).Value2 = 0.5432687471728113
$ws.Cells.Item(2, 7).Value2 = 16
$ws.Cells.Item(3, 6).Value2 = 0.6868845891952087
$ws.Cells.Item(4, 4).Value2 = 14
$ws.Cells.Item(4, 6).Value2 = 0.5253707450939883
$ws.Cells.Item(4, 7).Value2 = 13
$ws.Cells.Item(5, 3).Value2 = 2
$ws.Cells.Item(5, 5).Value2 = 0
$ws.Cells.Item(5, 6).Value2 = 0.5124433936071593
$ws.Cells.Item(5, 7).Value2 = 12
$ws.Cells.Item(6, 6).Value2 = 0.6385259457477233

# --- Sheet 4: Step3_DataPts_0.7 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 4).Value2 = 20
$ws.Cells.Item(2, 6).Value2 = 0.7339293761152323
$ws.Cells.Item(2, 7).Value2 = 19
$ws.Cells.Item(3, 4).Value2 = 20
$ws.Cells.Item(3, 6).Value2 = 0.7657704399026797
$ws.Cells.Item(3, 7).Value2 = 19
$ws.Cells.Item(4, 6).Value2 = 0.7006770111423088
$ws.Cells.Item(5, 3).Value2 = 2
$ws.Cells.Item(5, 5).Value2 = 0
$ws.Cells.Item(5, 6).Value2 = 0.7569312964202982
$ws.Cells.Item(5, 7).Value2 = 18
$ws.Cells.Item(6, 6).Value2 = 0.7503296723888881

# --- Sheet 5: Step3_DataPts_0.8 ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 6).Value2 = 0.8003170313702462
$ws.Cells.Item(3, 6).Value2 = 0.866907766882127
$ws.Cells.Item(4, 4).Value2 = 22
$ws.Cells.Item(4, 6).Value2 = 0.8700666786215903
$ws.Cells.Item(4, 7).Value2 = 21
$ws.Cells.Item(5, 3).Value2 = 2
$ws.Cells.Item(5, 5).Value2 = 0
$ws.Cells.Item(5, 6).Value2 = 0.8556685446971043
$ws.Cells.Item(5, 7).Value2 = 20
$ws.Cells.Item(6, 6).Value2 = 0.8669382243703267

# --- Sheet 6: Step3_DataPts_0.9 ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(3, 6).Value2 = 0.9935109594348895
$ws.Cells.Item(4, 4).Value2 = 24
$ws.Cells.Item(4, 6).Value2 = 0.9037175148254573
$ws.Cells.Item(4, 7).Value2 = 23
$ws.Cells.Item(5, 3).Value2 = 2
$ws.Cells.Item(5, 5).Value2 = 0
$ws.Cells.Item(5, 6).Value2 = 0.960520564543599
$ws.Cells.Item(5, 7).Value2 = 23
$ws.Cells.Item(6, 4).Value2 = 24
$ws.Cells.Item(6, 6).Value2 = 0.9061000636536229
$ws.Cells.Item(6, 7).Value2 = 23

Write-Output "Applied all pipeline value updates."
